$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.546.84'
$ws.Range('E2').Value = '  +3.44%  '
$ws.Range('D3').Value = '2.554.11'
$ws.Range('E3').Value = '  +3.72%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.38%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '543.69'
$ws.Range('E5').Value = '  +2.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.97'
$ws.Range('E6').Value = '  +2.36%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('E7').Value = '  -0.22%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.573'
$ws.Range('E8').Value = '  +1.25%  '
$ws.Range('D9').Value = '2.590.65'
$ws.Range('E9').Value = '  +4.42%  '
$ws.Range('E10').Value = '  +3.29%  '
$ws.Range('E12').Value = '  +0.36%  '
$ws.Range('E13').Value = '  +4.66%  '
$ws.Range('D14').Value = '3.003.59'
$ws.Range('E14').Value = '  +3.54%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '24.60'
$ws.Range('E15').Value = '  +4.01%  '
$ws.Range('D16').Value = '60.521.66'
$ws.Range('E16').Value = '  +3.50%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000145'
$ws.Range('E17').Value = '  +6.07%  '
$ws.Range('D18').Value = '2.557.34'
$ws.Range('E18').Value = '  +2.98%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.43'
$ws.Range('E19').Value = '  +1.77%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.39'
$ws.Range('E20').Value = '  +2.35%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '330.11'
$ws.Range('E21').Value = '  +2.60%  '
$ws.Range('B22').Value = 'Uniswap'
$ws.Range('C22').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.00'
$ws.Range('E22').Value = '  +5.02%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  +0.19%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '63.46'
$ws.Range('E24').Value = '  +4.77%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.445'
$ws.Range('E25').Value = '  +2.07%  '
$ws.Range('E26').Value = '  +5.13%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.992'
$ws.Range('E27').Value = '  -0.52%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.11'
$ws.Range('E28').Value = '  +5.84%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.18'
$ws.Range('E29').Value = '  +4.09%  '
$ws.Range('D30').Value = '0.0₃0816'
$ws.Range('E30').Value = '  +6.10%  '
$ws.Range('E31').Value = '  +2.39%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.23'
$ws.Range('E32').Value = '  +0.53%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '164.19'
$ws.Range('E33').Value = '  +4.14%  '
$ws.Range('E34').Value = '  +6.95%  '
$ws.Range('E35').Value = '  -0.11%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '18.92'
$ws.Range('E36').Value = '  +2.79%  '
$ws.Range('E37').Value = '  +3.48%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.67'
$ws.Range('E38').Value = '  +4.93%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.75'
$ws.Range('E39').Value = '  +0.44%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '310.14'
$ws.Range('E40').Value = '  +2.08%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '37.15'
$ws.Range('E41').Value = '  +1.95%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.852'
$ws.Range('E42').Value = '  +6.61%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.79'
$ws.Range('E43').Value = '  +2.85%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.614'
$ws.Range('E44').Value = '  +4.47%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.991'
$ws.Range('E45').Value = '  -0.47%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.87'
$ws.Range('E46').Value = '  +0.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '126.89'
$ws.Range('E47').Value = '  +2.49%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '19.19'
$ws.Range('E48').Value = '  +4.58%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0943'
$ws.Range('E49').Value = '  +2.32%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0528'
$ws.Range('E50').Value = '  +2.10%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0233'
$ws.Range('E51').Value = '  +2.48%  '
